$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix workout level for row 139 (Matt's Ride on 2024-07-01) -> Sauntering Hippo
$ws.Range("L139").Value = "Sauntering Hippo"

# Add two new rows of data (Steven - Walk - 2024-07-01)
$ws.Range("A141").Value = "Steven"
$ws.Range("B141").Value = 45474
$ws.Range("C141").Value = "Walk"
$ws.Range("D141").Value = 19
$ws.Range("E141").Value = 0.94
$ws.Range("F141").Value = 33
$ws.Range("G141").Value = 19
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = "Mighty Monkey"
$ws.Range("M141").Value = 4

$ws.Range("A142").Value = "Steven"
$ws.Range("B142").Value = 45474
$ws.Range("C142").Value = "Walk"
$ws.Range("D142").Value = 22
$ws.Range("E142").Value = 1.06
$ws.Range("F142").Value = 49
$ws.Range("G142").Value = 22
$ws.Range("H142").Value = 0
$ws.Range("I142").Value = 0
$ws.Range("J142").Value = 0
$ws.Range("K142").Value = 0
$ws.Range("L142").Value = "Mighty Monkey"
$ws.Range("M142").Value = 4

# Copy the date formatting from an existing date cell so no new number
# format gets minted (keeps the same style id as the rest of column B)
$ws.Range("B140").Copy()
$ws.Range("B141:B142").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selection to match the post-edit state
$ws.Range("J148").Select()
